# Chapter 5 rewrite of the Opera load-time table (rows 31-42): the old
# "Dev vs Vulcanized x Time-to-load/Time-to-display" grid becomes a per
# browser Initial -> Complete table with a computed Rate-of-Reduction
# column. New text is entered in the order it first appears so shared
# strings intern in the same order as the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Opera section: row 33 "Initiatl" / row 34 "Complete" ---
$ws.Range("A33").Value = "Initiatl"
$ws.Range("B33").Value = 2800
$ws.Range("C33").Value = 1100
$ws.Range("D33").Formula = "=(C33-B33)/B33"

$ws.Range("A34").Value = "Complete"
$ws.Range("B34").Value = 7000
$ws.Range("C34").Value = 2800
$ws.Range("D34").Formula = "=(C34-B34)/B34"

# --- Google Chrome section: row 37 "Initial" / row 38 "Complete" ---
$ws.Range("A37").Value = "Initial"
$ws.Range("B37").Value = 2200
$ws.Range("C37").Value = 1150
$ws.Range("D37").Formula = "=(C37-B37)/B37"

$ws.Range("A38").Value = "Complete"
$ws.Range("B38").Value = 4200
$ws.Range("C38").Value = 2900
$ws.Range("D38").Formula = "=(C38-B38)/B38"

# --- Firefox section: row 41 "Initial" / row 42 "Complete" ---
$ws.Range("A41").Value = "Initial"
$ws.Range("B41").Value = 7400
$ws.Range("C41").Value = 6050
$ws.Range("D41").Formula = "=(C41-B41)/B41"

$ws.Range("A42").Value = "Complete"
$ws.Range("B42").Value = 12000
$ws.Range("C42").Value = 10500
$ws.Range("D42").Formula = "=(C42-B42)/B42"

# --- Row 32 header (reuses Dev / Vulcanized, adds "Rate of Reduction") ---
$ws.Range("B32").Value = "Dev"
$ws.Range("C32").Value = "Vulcanized"
$ws.Range("D32").Value = "Rate of Reduction"

# --- Browser section headers (reuse existing strings) ---
$ws.Range("A36").Value = "Google Chrome"
$ws.Range("A40").Value = "Firefox"

# Percent-format the ratio column, including the blank spacer rows that
# still carry the formatting (35, 36, 39, 40).
$ws.Range("D33:D34").NumberFormat = "0.00%"
$ws.Range("D35").NumberFormat = "0.00%"
$ws.Range("D36").NumberFormat = "0.00%"
$ws.Range("D37:D38").NumberFormat = "0.00%"
$ws.Range("D39").NumberFormat = "0.00%"
$ws.Range("D40").NumberFormat = "0.00%"
$ws.Range("D41:D42").NumberFormat = "0.00%"

# Match the author's final selection and page orientation.
$ws.Range("D37").Select() | Out-Null
$ws.PageSetup.Orientation = 1
